$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (header "Förändrad") holds a date serial number that was bumped
# by one day (46075 -> 46076) for every data row (rows 2 through 499).
for ($r = 2; $r -le 499; $r++) {
    $ws.Cells.Item($r, 3).Value = 46076
}
